$wb = $excel.ActiveWorkbook

# --- Update the Date property on the "Metadata" sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Add a new mapping column on the "Elements" sheet ---
$ws = $wb.Worksheets.Item("Elements")

# Copy the header style/value from the neighboring header cell (AK1) and data
# cells from AK2:AK6 so the new column matches the existing "empty" placeholder
# pattern (these are shared-string empty cells, not truly blank cells).
$ws.Range("AK1:AK6").Copy()
$ws.Range("AL1:AL6").PasteSpecial(-4122)

# New header text in column AL (column 38)
$ws.Cells.Item(1, 38).Value = "Mapping: Spécification métier vers l'extension ROR Comment"

# Data rows: only the last row (Extension.value[x]) has a mapping value
$ws.Cells.Item(6, 38).Value = "commentaire"

# Column width to match the diff (bestFit width ~64.89)
$ws.Columns.Item(38).ColumnWidth = 64
